$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("D3").Value = -7.865
$ws.Range("C12").Value = -10.896
$ws.Range("D14").Value = -7.722000000000001
$ws.Range("D26").Value = -8.026999999999999
$ws.Range("C27").Value = -12.694
$ws.Range("D31").Value = -8.397
$ws.Range("C32").Value = -13.676
$ws.Range("D35").Value = -7.935
$ws.Range("C36").Value = -12.732
$ws.Range("D37").Value = -7.712999999999999
$ws.Range("C38").Value = -12.703
$ws.Range("D45").Value = -7.57
$ws.Range("C46").Value = -13.731
$ws.Range("D52").Value = -7.280999999999999
$ws.Range("C54").Value = -13.137
$ws.Range("C55").Value = -13.391
$ws.Range("C56").Value = -13.228
$ws.Range("D57").Value = -8.184000000000001
$ws.Range("C67").Value = -11.705
$ws.Range("C69").Value = -11.038
$ws.Range("C72").Value = -11.555
$ws.Range("D81").Value = -7.221000000000001
$ws.Range("C83").Value = -13.392
$ws.Range("D83").Value = -8.461
$ws.Range("C86").Value = -13.846
$ws.Range("C91").Value = -10.885
$ws.Range("C93").Value = -11.979
$ws.Range("C99").Value = -12.635
$ws.Range("D100").Value = -8.234999999999999
$ws.Range("D102").Value = -7.790000000000001
